$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 86: copy formatting from the matching columns above so the
# new cells pick up the same styles used elsewhere in the sheet
# (column A uses the "key" style, column B here uses the alternate
# "value" style, matching rows 79/83).
$ws.Range("A85").Copy()
$ws.Range("A86").PasteSpecial(-4122)

$ws.Range("B79").Copy()
$ws.Range("B86").PasteSpecial(-4122)

$ws.Range("A86").Value = "AbsenceLoginInfo"
$ws.Range("B86").Value = "Devamsizlik SMS gonderimi icin e-okul oturumunuzu acmaniz gerekmektedir."
